$d = $word.ActiveDocument

# --- Step 1: split the existing "hola" run into two runs: "H" and "ola" ---
# (same text/formatting, just emitted as two <w:r> elements instead of one)
$splitXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>H</w:t></w:r><w:r><w:t>ola</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$firstPara = $d.Paragraphs(1)
# Use the paragraph's text only (exclude the trailing paragraph mark) so the
# mark-bearing <w:p> itself is reused rather than duplicated.
$firstRange = $d.Range($firstPara.Range.Start, $firstPara.Range.End - 1)
$firstRange.InsertXML($splitXml)

# --- Step 2: append an empty paragraph, a "Holaa" paragraph (flagged by the
#     spell checker via proofErr spellStart/spellEnd), and a trailing empty
#     paragraph after the first paragraph ---
$tailXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Holaa</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$endPos = $d.Content.End
$endRange = $d.Range($endPos, $endPos)
$endRange.InsertXML($tailXml)
